$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: mark the D-column cells whose new values look like plain numbers
# as Text format first, so Excel stores the literal string (matching the
# source data which keeps trailing zeros / exact decimal formatting) rather
# than silently parsing them into a Double and losing precision/formatting.
$riskyCells = @("D5", "D6", "D8", "D12", "D13", "D17", "D19", "D21", "D22", "D24", "D25", "D26", "D27", "D29", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $riskyCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Step 2: write the new cell values row by row
$ws.Range("D2").Value = "68.321.93"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "2.452.13"
$ws.Range("E3").Value = "  -2.04%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "560.23"
$ws.Range("E5").Value = "  -2.77%  "
$ws.Range("D6").Value = "163.16"
$ws.Range("E6").Value = "  -2.24%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.503"
$ws.Range("E8").Value = "  -1.92%  "
$ws.Range("D9").Value = "2.451.91"
$ws.Range("E9").Value = "  -2.00%  "
$ws.Range("E10").Value = "  -5.79%  "
$ws.Range("E11").Value = "  -2.01%  "
$ws.Range("D12").Value = "0.339"
$ws.Range("E12").Value = "  -5.34%  "
$ws.Range("D13").Value = "4.82"
$ws.Range("E13").Value = "  -2.53%  "
$ws.Range("D14").Value = "2.909.41"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "68.353.41"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("E16").Value = "  -3.77%  "
$ws.Range("D17").Value = "23.30"
$ws.Range("E17").Value = "  -5.80%  "
$ws.Range("D18").Value = "2.463.62"
$ws.Range("E18").Value = "  -1.62%  "
$ws.Range("D19").Value = "10.97"
$ws.Range("E19").Value = "  -2.41%  "
$ws.Range("E20").Value = "  -4.20%  "
$ws.Range("D21").Value = "342.20"
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("D22").Value = "3.79"
$ws.Range("E22").Value = "  -3.27%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "1.86"
$ws.Range("E24").Value = "  -3.97%  "
$ws.Range("D25").Value = "67.82"
$ws.Range("E25").Value = "  -3.90%  "
$ws.Range("D26").Value = "1.06"
$ws.Range("E26").Value = "  +4.42%  "
$ws.Range("D27").Value = "3.71"
$ws.Range("E27").Value = "  -6.20%  "
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("D29").Value = "8.16"
$ws.Range("E29").Value = "  -7.07%  "
$ws.Range("D30").Value = "0.0₃0834"
$ws.Range("E30").Value = "  -6.59%  "
$ws.Range("D31").Value = "7.24"
$ws.Range("E31").Value = "  -7.48%  "
$ws.Range("D32").Value = "3.33"
$ws.Range("E32").Value = "  +125.28%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").Value = "432.47"
$ws.Range("E34").Value = "  -5.36%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "1.17"
$ws.Range("E35").Value = "  -3.46%  "
$ws.Range("E36").Value = "  -3.40%  "
$ws.Range("D37").Value = "156.75"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").Value = "19.00"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.110"
$ws.Range("E40").Value = "  -5.52%  "
$ws.Range("D41").Value = "17.87"
$ws.Range("E41").Value = "  -3.08%  "
$ws.Range("D42").Value = "0.306"
$ws.Range("E42").Value = "  -3.45%  "
$ws.Range("D43").Value = "4.45"
$ws.Range("E43").Value = "  -5.07%  "
$ws.Range("D44").Value = "1.51"
$ws.Range("E44").Value = "  -5.25%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "2.06"
$ws.Range("E46").Value = "  -6.46%  "
$ws.Range("D47").Value = "133.93"
$ws.Range("E47").Value = "  -5.66%  "
$ws.Range("D48").Value = "3.35"
$ws.Range("E48").Value = "  -3.57%  "
$ws.Range("D50").Value = "0.484"
$ws.Range("E50").Value = "  -6.96%  "
$ws.Range("D51").Value = "0.560"
$ws.Range("E51").Value = "  -3.16%  "

# Step 3: restore the Normal style on the cells we forced to Text format,
# so number formatting/style indices go back to the workbook default
# (values remain text; this only resets the display style).
foreach ($addr in $riskyCells) {
    $ws.Range($addr).Style = "Normal"
}